# ISIS-2181: upgrades cucumber (specsupport); also updates mml integtests
# for simpleapp. Repositions several shapes/connectors on the pom-hierarchy
# diagram and adds a duplicate "Convenience pom" callout textbox.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU -> Points helper (PowerPoint COM geometry properties are expressed
# in points; 1 point = 12700 EMU).
function EMU([double]$emu) {
    return $emu / 12700.0
}

function MoveShape($index, $x, $y) {
    $sh = $s.Shapes.Item($index)
    $sh.Left = EMU $x
    $sh.Top  = EMU $y
}

function MoveAndResizeShape($index, $x, $y, $cx, $cy) {
    $sh = $s.Shapes.Item($index)
    $sh.Left   = EMU $x
    $sh.Top    = EMU $y
    $sh.Width  = EMU $cx
    $sh.Height = EMU $cy
}

# 1) Rectangle 91 (empty placeholder rectangle) - move only
MoveShape 1 8744064 3279262

# 2) Rectangle 9 ("org.apache.isis.mavendeps / isis-mavendeps") - move only
MoveShape 10 8415362 4145164

# 3) Rectangle 12 ("org.apache.isis.examples.apps / helloworld") - move only
MoveShape 13 8669702 3203416

# 4) Connector: Curved 25 - move + resize
MoveAndResizeShape 20 6669497 1057889 456234 5718315

# 5) Connector: Curved 77 - move + resize
MoveAndResizeShape 33 8384274 3743621 361348 2383646

# 6) Connector: Curved 80 - move + resize + flip/rotation change
#    before: rot=270 (16200000/60000), flipV only
#    after:  rot=90  (5400000/60000), flipH + flipV
$connector80 = $s.Shapes.Item(34)
$connector80.Rotation = 90
$connector80.HorizontalFlip = -1
$connector80.VerticalFlip = -1
MoveAndResizeShape 34 9324553 4683899 361347 503090

# 7) Connector: Curved 83 - move + resize
MoveAndResizeShape 35 10246802 4264740 361347 1341408

# 8) Connector: Curved 92 - move + resize (rotation unchanged at 180)
MoveAndResizeShape 36 7403116 2708393 1266586 799826

# 9) Add a new "Convenience pom / scope=import, / type=pom" textbox; it is
#    a duplicate of the existing TextBox 56 callout placed alongside it.
$source = $s.Shapes.Item(53)
$clone = $source.Duplicate()
$clone.Name = "TextBox 54"
$clone.Left   = EMU 5573144
$clone.Top    = EMU 4112289
$clone.Width  = EMU 2583807
$clone.Height = EMU 707886
